$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (origen) was curated from a "dimension" (skos:Concept, with an
# external mapping-origen.xlsx lookup) into a plain "measure" of type xsd:int.
$ws.Range("D2").Value = "iaest-measure:origen"
$ws.Range("D3").Value = "medida"
$ws.Range("D4").Value = "xsd:int"
$ws.Range("D5").Clear()

# Column I (tipo-de-alojamiento) received the same curation treatment.
$ws.Range("I2").Value = "iaest-measure:tipo-de-alojamiento"
$ws.Range("I3").Value = "medida"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("I5").Clear()

# Column J (municipio-nombre) moves the other way: it becomes a curated
# dimension resolved against sdmx-dimension:refArea / URI-Municipio, mirroring
# columns K (provincia-nombre / URI-Provincia) and M (comarca-nombre / URI-comarca).
$ws.Range("J2").Value = "sdmx-dimension:refArea"
$ws.Range("J3").Value = "dim"
$ws.Range("J4").Value = "URI-Municipio"
